# Updated MCH102 to MCH251
#
# Adds 8 archival-collection records (MCH203-1 .. MCH203-8) to Sheet1,
# below the existing header row.
#
# Columns (per row 1 headers):
#   A identifier            C title                E levelOfDescription
#   B alternativeIdentifiers D date_s               F extentAndMedium
#                                                    G notes
#                                                    H file_path (left blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row data. B and H are intentionally left empty, matching the source.
# ---------------------------------------------------------------------------

$rows = @(
    @{ Row=2; Id="MCH203-1"; Title="";                                          Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24B | GRAP COUNT NUMER: NONE" },
    @{ Row=3; Id="MCH203-2"; Title="";                                          Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24B | GRAP COUNT NUMER: NONE" },
    @{ Row=4; Id="MCH203-3"; Title="BOYKOT- RUNDBRIEF NO. 2-31";                Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24B | GRAP COUNT NUMER: NONE" },
    @{ Row=5; Id="MCH203-4"; Title="BOYKOT- RUNDBRIEF NO. 32-53, INFORMATIONEN- EVANGELISCHEN RELIGIONSUNTERRICHT IN BERLIN 1980, 86, 88, MODELHE FURDEN RELIGIONS, UNTERRICHT 4"; Date="1980"; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24C | GRAP COUNT NUMER: NONE" },
    @{ Row=6; Id="MCH203-5"; Title="MAYIBUYE CENTRE UWC, ABS: NINIKRAATZ ROONSTR. 10A H163 BERLIN GERMANY"; Date=""; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24C | GRAP COUNT NUMER: NONE" },
    @{ Row=7; Id="MCH203-6"; Title="ARTIFACTS TO BE CATALOUGED";                Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24C | GRAP COUNT NUMER: NONE" },
    @{ Row=8; Id="MCH203-7"; Title="ARTIFACTS TO BE CATALOUGED";                Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24C | GRAP COUNT NUMER: NONE" },
    @{ Row=9; Id="MCH203-8"; Title="";                                          Date=""    ; Level="Series"; Extent="1 Box"; Notes="LOCATION: 24C | GRAP COUNT NUMER: NONE" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id      # A identifier
    if ($r.Title -ne "") {
        $ws.Cells.Item($r.Row, 3).Value = $r.Title   # C title
    }
    if ($r.Date -ne "") {
        # Force text storage so a purely-numeric date (e.g. "1980") isn't
        # silently coerced into a Number cell - the source keeps it as text.
        $ws.Cells.Item($r.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($r.Row, 4).Value = $r.Date    # D date_s
        $ws.Cells.Item($r.Row, 4).NumberFormat = "General"
    }
    $ws.Cells.Item($r.Row, 5).Value = $r.Level   # E levelOfDescription
    $ws.Cells.Item($r.Row, 6).Value = $r.Extent  # F extentAndMedium
    $ws.Cells.Item($r.Row, 7).Value = $r.Notes   # G notes
}

# ---------------------------------------------------------------------------
# 2. Formatting - the new rows use the workbook's plain data-row look:
#    Calibri 10, automatic/theme text colour, no fill. Build the style once
#    on A2 and fan it out with a format-only paste so every new cell shares
#    a single style record instead of minting one per cell.
# ---------------------------------------------------------------------------

$ws.Cells.Item(2, 1).Font.Name = "Calibri"
$ws.Cells.Item(2, 1).Font.Size = 10
$ws.Cells.Item(2, 1).Font.ThemeColor = 1
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)
$ws.Range("C2:D9").PasteSpecial(-4122)
$ws.Range("E2:E9").PasteSpecial(-4122)
$ws.Range("F2:F9").PasteSpecial(-4122)
$ws.Range("G2:H9").PasteSpecial(-4122)

$excel.CutCopyMode = 0
